# Pull/refresh the DataSource sheet with the new "preproduccion" Gestion
# Documental environment (replacing the old Oracle dev gateway) and the
# matching account number.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column B (URL) is set first, then column A (Ambiente/hostname), so the two
# brand-new shared strings land in the same order the source workbook has
# them: the full URL, then the bare hostname.
$ws.Range("B2").Value = "https://i-preproducciongestion.segurossura.com.ar/pc/PolicyCenter.do"
$ws.Range("A2").Value = "i-preproducciongestion.segurossura.com.ar"

# NroCuenta (account number) for the new environment.
$ws.Range("E2").Value = 2240451788

# The URL text no longer fits the default column width, so column B was
# widened by hand.
$ws.Columns.Item(2).ColumnWidth = 40.5

# Scroll the view over towards the right-hand columns and leave the
# selection on R2, as last left by the editor.
$excel.ActiveWindow.ScrollColumn = 9
$ws.Range("R2").Select() | Out-Null
